$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve values as literal text (matching source data which is inlineStr),
# since these columns hold formatted price/percentage strings, not live numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.92%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.16%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.037"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07908"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.53%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.850"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.97%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.102"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.07%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.785"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.18%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9201"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1347"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1889"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.89%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09058"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.23%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03469"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.69%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09806"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001397"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.62%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006141"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "5.33%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.723"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.53%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "9.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3438"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1342"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.162"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.74%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2192"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-9.03%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04399"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.50%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.65%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004612"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.53%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.84%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004441"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.03%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05085"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.10%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007635"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.22%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01010"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-7.92%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1341"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002171"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.33%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01019"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006150"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.67%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.69%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001659"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.23%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
